$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: shared string text changed from "RFin DSTsel Yout END" to "RFin DSTsel Yout"
$ws.Range("F6").Value = "RFin DSTsel Yout"

# Row 7: new microcode line "MOV" / addressing mode 0 / time-step 5 / "TCend"
$ws.Range("A7").Value = "MOV"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 5
$ws.Range("F7").Value = "TCend"

# B5 loses its explicit alignment/protection override, reverting to the base style
$ws.Range("B5").NumberFormat = "General"

# Move the active selection to E20
$ws.Range("E20").Select()
